$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 0.0021097888071897361
$ws.Range("B1").Value = 0.19969781443587659
$ws.Range("C1").Value = 0.028866838732877646
$ws.Range("D1").Value = 818877554.49260378
$ws.Range("E1").Value = 23.824832878077498
$ws.Range("F1").Value = 61.845045926722669
$ws.Range("G1").Value = 0.68957277751298429
$ws.Range("H1").Value = 0.070567408858400738
$ws.Range("I1").Value = 0.14778064437536209
$ws.Range("J1").Value = 0.00000018567800421632464
$ws.Range("K1").Value = 0.99150332624430371
$ws.Range("L1").Value = 0.98662767211449409
$ws.Range("M1").Value = 0.97269733624408994
$ws.Range("N1").Value = 0.091715358681844147
$ws.Range("O1").Value = 0.11505913525944325
$ws.Range("P1").Value = 0.16440692539656282

$ws.Range("A2").Value = 0.0023999735378474928
$ws.Range("B2").Value = 0.30674523852631458
$ws.Range("C2").Value = 0.026315789624693837
$ws.Range("D2").Value = 59991281.952788375
$ws.Range("E2").Value = 20.000000000765493
$ws.Range("F2").Value = 62.77548464302069
$ws.Range("G2").Value = 0.046762112047392383
$ws.Range("H2").Value = 0.000000012219467313599025
$ws.Range("I2").Value = 0.24027361900235344
$ws.Range("J2").Value = 0.099999999999637254
$ws.Range("K2").Value = 0.99626937496544443
$ws.Range("L2").Value = 0.9587232081938093
$ws.Range("M2").Value = 0.9820650985548256
$ws.Range("N2").Value = 0.060772681232689078
$ws.Range("O2").Value = 0.20214851938149031
$ws.Range("P2").Value = 0.13324996221659002

$ws.Range("A3").Value = 0.0025203020454935741
$ws.Range("B3").Value = 0.06908592067598554
$ws.Range("C3").Value = 0.099510365861983335
$ws.Range("D3").Value = 1041012019.957088
$ws.Range("E3").Value = 20.000000000227651
$ws.Range("F3").Value = 67.826297979498108
$ws.Range("G3").Value = 0.60732182036244786
$ws.Range("H3").Value = 0.013304101090591923
$ws.Range("I3").Value = 0.0000010750661675905394
$ws.Range("J3").Value = 0.00000000000068956276972077998
$ws.Range("K3").Value = 0.99144688812041049
$ws.Range("L3").Value = 0.98505815850502187
$ws.Range("M3").Value = 0.99343822590965092
$ws.Range("N3").Value = 0.092019458598676815
$ws.Range("O3").Value = 0.1216241056700041
$ws.Range("P3").Value = 0.0805987366491909

$ws.Range("A4").Value = 0.0064071421325667833
$ws.Range("B4").Value = 0.25980706426117811
$ws.Range("C4").Value = 0.026315789480790875
$ws.Range("D4").Value = 86708347.125568867
$ws.Range("E4").Value = 30.918362750177526
$ws.Range("F4").Value = 58.004475857179656
$ws.Range("G4").Value = 0.000000010000022218739533
$ws.Range("H4").Value = 0.17453852000111
$ws.Range("I4").Value = 0.00000000046450674282675123
$ws.Range("J4").Value = 0.099999999917342125
$ws.Range("K4").Value = 0.9935565156363213
$ws.Range("L4").Value = 0.94312540558446789
$ws.Range("M4").Value = 0.93866074510170361
$ws.Range("N4").Value = 0.079868952165669124
$ws.Range("O4").Value = 0.23728853421810506
$ws.Range("P4").Value = 0.24642618032448055

$ws.Range("A5").Value = 0.018201868192476232
$ws.Range("B5").Value = 1.9708590250213547
$ws.Range("C5").Value = 0.075990278067578604
$ws.Range("D5").Value = 1103118534.2265716
$ws.Range("E5").Value = 20.001289286491811
$ws.Range("F5").Value = 62.022909125504256
$ws.Range("G5").Value = 0.62163832419773746
$ws.Range("H5").Value = 1.4687468885933381
$ws.Range("I5").Value = 0.00042703010765250308
$ws.Range("J5").Value = 0.0017244618951913167
$ws.Range("K5").Value = 0.97307379874000466
$ws.Range("L5").Value = 0.97032912497335466
$ws.Range("M5").Value = 0.99000839728903078
$ws.Range("N5").Value = 0.16326952945174861
$ws.Range("O5").Value = 0.17138893277098996
$ws.Range("P5").Value = 0.099456958951395266

$ws.Range("A6").Value = 0.010570746461536065
$ws.Range("B6").Value = 0.29472793598232966
$ws.Range("C6").Value = 0.026316402995454225
$ws.Range("D6").Value = 98228881.944859281
$ws.Range("E6").Value = 23.123651898412643
$ws.Range("F6").Value = 54.448437247405032
$ws.Range("G6").Value = 0.000000023180235490698976
$ws.Range("H6").Value = 0.11652082502789272
$ws.Range("I6").Value = 0.0000000039204152073294536
$ws.Range("J6").Value = 0.099999889435985811
$ws.Range("K6").Value = 0.99563208717212115
$ws.Range("L6").Value = 0.98880646029675801
$ws.Range("M6").Value = 0.99023107381598163
$ws.Range("N6").Value = 0.065758905857686359
$ws.Range("O6").Value = 0.10526919922849959
$ws.Range("P6").Value = 0.098342447204542494

$ws.Range("A7").Value = 0.0087969082514614411
$ws.Range("B7").Value = 0.05809550120867029
$ws.Range("C7").Value = 0.072854988985842833
$ws.Range("D7").Value = 869497827.19434929
$ws.Range("E7").Value = 39.086659225568162
$ws.Range("F7").Value = 54.557059569815017
$ws.Range("G7").Value = 0.80314177116212115
$ws.Range("H7").Value = 0.00011943425624645336
$ws.Range("I7").Value = 0.055637136089204252
$ws.Range("J7").Value = 0.00000000043013514913368757
$ws.Range("K7").Value = 0.99705613154973849
$ws.Range("L7").Value = 0.9961715874489161
$ws.Range("M7").Value = 0.99672869148871157
$ws.Range("N7").Value = 0.053985458836235956
$ws.Range("O7").Value = 0.061564018919926543
$ws.Range("P7").Value = 0.056908658622177449

$ws.Range("A8").Value = 0.010293908569875911
$ws.Range("B8").Value = 0.062545996581536478
$ws.Range("C8").Value = 0.099999997765947368
$ws.Range("D8").Value = 1105943839.9989924
$ws.Range("E8").Value = 34.224737839072297
$ws.Range("F8").Value = 62.585372864441652
$ws.Range("G8").Value = 0.66739486173893814
$ws.Range("H8").Value = 0.078521634261090573
$ws.Range("I8").Value = 0.00000080980355135467739
$ws.Range("J8").Value = 0.000000000058042029473927693
$ws.Range("K8").Value = 0.99748717120454367
$ws.Range("L8").Value = 0.99501463086482556
$ws.Range("M8").Value = 0.99238563291700888
$ws.Range("N8").Value = 0.04987934745395102
$ws.Range("O8").Value = 0.070256736650140608
$ws.Range("P8").Value = 0.086827284363224247

$ws.Range("A9").Value = 0.025199604094165075
$ws.Range("B9").Value = 3.2940758033856365
$ws.Range("C9").Value = 0.052590670563293862
$ws.Range("D9").Value = 488741322.03084129
$ws.Range("E9").Value = 29.432058148267789
$ws.Range("F9").Value = 66.759459284471447
$ws.Range("G9").Value = 0.058078974343964036
$ws.Range("H9").Value = 9.6064519794520962
$ws.Range("I9").Value = 0.0011616304467309338
$ws.Range("J9").Value = 0.085130223140294115
$ws.Range("K9").Value = 0.99576792421273841
$ws.Range("L9").Value = 0.93831158973104345
$ws.Range("M9").Value = 0.97353302689721166
$ws.Range("N9").Value = 0.064728317059761342
$ws.Range("O9").Value = 0.24712653877369592
$ws.Range("P9").Value = 0.16187125554514148

$ws.Range("A10").Value = 0.0043802661548073187
$ws.Range("B10").Value = 0.44237047636065913
$ws.Range("C10").Value = 0.087083200564253657
$ws.Range("D10").Value = 63180115.750805609
$ws.Range("E10").Value = 37.459781793540159
$ws.Range("F10").Value = 40.161666358156182
$ws.Range("G10").Value = 0.030267711350570131
$ws.Range("H10").Value = 0.43153982247464184
$ws.Range("I10").Value = 0.0046296795348582222
$ws.Range("J10").Value = 0.000010495607472909202
$ws.Range("K10").Value = 0.99450862365726522
$ws.Range("L10").Value = 0.98651909974957253
$ws.Range("M10").Value = 0.98398951633829879
$ws.Range("N10").Value = 0.073732371312113965
$ws.Range("O10").Value = 0.11552528402009304
$ws.Range("P10").Value = 0.12589828761775976

$ws.Range("A11").Value = 0.008827856827243824
$ws.Range("B11").Value = 0.0755417192547954
$ws.Range("C11").Value = 0.045248107026542507
$ws.Range("D11").Value = 110535998
$ws.Range("E11").Value = 20.000000000000021
$ws.Range("F11").Value = 70.574558882166301
$ws.Range("G11").Value = 0.62548906338310539
$ws.Range("H11").Value = 0.097814635328407312
$ws.Range("I11").Value = 0.000028165400896217941
$ws.Range("J11").Value = 0.000000000000022690196750526574
$ws.Range("K11").Value = 0.98339886367076657
$ws.Range("L11").Value = 0.98376764779690362
$ws.Range("M11").Value = 0.97963641930253809
$ws.Range("N11").Value = 0.12819955134843913
$ws.Range("O11").Value = 0.12676761684699067
$ws.Range("P11").Value = 0.14198572072742852

$ws.Range("A12").Value = 0.06632247829524783
$ws.Range("B12").Value = 1.7137639329917644
$ws.Range("C12").Value = 0.069866626416853858
$ws.Range("D12").Value = 1451042807.7431676
$ws.Range("E12").Value = 20.000023683794321
$ws.Range("F12").Value = 68.607190433600962
$ws.Range("G12").Value = 0.50794113155111209
$ws.Range("H12").Value = 3.3604547861156884
$ws.Range("I12").Value = 0.0015720979459276606
$ws.Range("J12").Value = 0.00000011075155570108483
$ws.Range("K12").Value = 0.91251201628020961
$ws.Range("L12").Value = 0.90445659926885413
$ws.Range("M12").Value = 0.97619091556320547
$ws.Range("N12").Value = 0.29430104295192794
$ws.Range("O12").Value = 0.30755156758474572
$ws.Range("P12").Value = 0.15352847811538617

$ws.Range("A13").Value = 0.006009484252225218
$ws.Range("B13").Value = 0.023712350592520037
$ws.Range("C13").Value = 0.099999999886581772
$ws.Range("D13").Value = 217678933.50945121
$ws.Range("E13").Value = 28.447336580513571
$ws.Range("F13").Value = 54.924844146080751
$ws.Range("G13").Value = 0.0004388135617570779
$ws.Range("H13").Value = 0.0080096388493488026
$ws.Range("I13").Value = 0.0000000017379072526304474
$ws.Range("J13").Value = 0.035025895502838085
$ws.Range("K13").Value = 0.99860082644231662
$ws.Range("L13").Value = 0.99561191590596565
$ws.Range("M13").Value = 0.99032699534640034
$ws.Range("N13").Value = 0.037218030873577709
$ws.Range("O13").Value = 0.065910570116590336
$ws.Range("P13").Value = 0.097858441675021793

$ws.Range("A14").Value = 0.0038862277447454445
$ws.Range("B14").Value = 0.02725596656686213
$ws.Range("C14").Value = 0.042301835077202837
$ws.Range("D14").Value = 2505889729.8533902
$ws.Range("E14").Value = 33.024138056072623
$ws.Range("F14").Value = 59.464000620651539
$ws.Range("G14").Value = 0.73479274856265109
$ws.Range("H14").Value = 0.00000035911675585669334
$ws.Range("I14").Value = 0.030393711958913223
$ws.Range("J14").Value = 0.00000000013575856999693136
$ws.Range("K14").Value = 0.99716413806613691
$ws.Range("L14").Value = 0.99370151489655978
$ws.Range("M14").Value = 0.99165365355522117
$ws.Range("N14").Value = 0.052985878444397248
$ws.Range("O14").Value = 0.078965183798974159
$ws.Range("P14").Value = 0.090900401431077332

$ws.Range("A15").Value = 0.0068560554960869612
$ws.Range("B15").Value = 0.041542035450387677
$ws.Range("C15").Value = 0.032824140711530357
$ws.Range("D15").Value = 367617109.83428496
$ws.Range("E15").Value = 20.000000000005215
$ws.Range("F15").Value = 69.722828941943575
$ws.Range("G15").Value = 0.051486536699527199
$ws.Range("H15").Value = 0.00039011309118147674
$ws.Range("I15").Value = 0.00000000037299730241668449
$ws.Range("J15").Value = 0.099999999999977676
$ws.Range("K15").Value = 0.99820817721960153
$ws.Range("L15").Value = 0.99249749898466444
$ws.Range("M15").Value = 0.99383472763390213
$ws.Range("N15").Value = 0.042117746290541561
$ws.Range("O15").Value = 0.086182805739788912
$ws.Range("P15").Value = 0.07812566570875995

$ws.Range("A16").Value = 0.0052688700173783897
$ws.Range("B16").Value = 0.61214822745053432
$ws.Range("C16").Value = 0.075473443836953563
$ws.Range("D16").Value = 1519978906.0318289
$ws.Range("E16").Value = 42.233837038396345
$ws.Range("F16").Value = 67.258519619581634
$ws.Range("G16").Value = 0.65344582250337191
$ws.Range("H16").Value = 3.1034722121742919
$ws.Range("I16").Value = 0.0027720831879434263
$ws.Range("J16").Value = 0.000023605345234617402
$ws.Range("K16").Value = 0.99372489602609093
$ws.Range("L16").Value = 0.96468552283099351
$ws.Range("M16").Value = 0.91289813590883506
$ws.Range("N16").Value = 0.078818480917675668
$ws.Range("O16").Value = 0.18697949726458349
$ws.Range("P16").Value = 0.29365089042986631
